# Audi.xlsx: fix the registration-number casing typo on Sheet1 and move the
# active selection, as captured by the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# A2 held "kX67 LFA" - correct the leading letter's casing to "KX67 LFA".
$ws.Range("A2").Value = "KX67 LFA"

# Move/record the active selection to B4 (was H3).
$ws.Range("B4").Select()
